# SpawnList update ("EA 23.222 Patch 1" / "nyan nyan" commit):
#  - bump the shared "last seen" version label EA 23.209 Patch 2 -> EA 23.222 Patch 1
#  - add two new chara categories, both last touched in EA 23.222 Patch 1:
#      c_dungeon_forest
#      c_machine
#    inserted right after the existing "c_dungeon" row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 is a blank header/filter row, data starts at row 3.
# "c_dungeon" currently sits on row 12, "c_sand" follows on row 13.
# Insert two fresh rows at row 13 so the new entries land right after c_dungeon.
$ws.Rows.Item(13).Insert()
$ws.Rows.Item(13).Insert()

$ws.Cells.Item(13, 1).Value = "c_dungeon_forest"
$ws.Cells.Item(13, 2).Value = "EA 23.222 Patch 1"

$ws.Cells.Item(14, 1).Value = "c_machine"
$ws.Cells.Item(14, 2).Value = "EA 23.222 Patch 1"
